# Replace the accelerometer sample data (rows 2-21) with a new, longer
# window of samples (rows 2-31). Header row (x, y, z) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-4.746825218200684, -21.48864555358887, -7.668013572692871),
    @(9.662276268005373, -37.98822784423828, -8.327095985412598),
    @(1.045047998428345, -10.61942100524902, 2.765533447265625),
    @(5.131261825561523, -28.12363815307617, 22.79096603393555),
    @(-30.11330032348633, -17.59586143493652, -19.26617050170898),
    @(1.400394916534424, -9.052268028259276, -1.365690350532532),
    @(29.48022842407227, -27.68916893005371, -8.951043128967285),
    @(20.90522193908692, -18.83166885375977, 3.978492736816406),
    @(6.43248987197876, -11.59229469299316, -0.6353058815002441),
    @(3.05394172668457, -49.75492095947266, 13.27557945251465),
    @(-3.964067220687866, -5.915932655334473, -9.63399600982666),
    @(-17.45916557312012, -25.58492851257324, -11.17863464355469),
    @(-3.193105697631836, -11.61043167114258, 9.58786392211914),
    @(-3.173403739929199, -4.077390670776367, 2.228257656097412),
    @(-78.2386703491211, -53.79793548583984, -33.76652908325195),
    @(4.821199893951416, -2.601359367370605, -5.058528423309326),
    @(-9.180764198303224, -23.62848663330078, -4.022332191467285),
    @(1.521630764007568, -8.447349548339844, 13.54604339599609),
    @(16.57039260864258, -22.10472106933594, 21.93498611450196),
    @(-78.08035278320312, -46.46374893188477, -22.17394256591797),
    @(-11.79047203063965, -6.283020973205566, -7.190555095672607),
    @(-2.819984912872314, -18.76873397827148, -7.554898738861084),
    @(5.317728996276856, -1.513343572616577, -2.32539701461792),
    @(5.113605499267578, -16.18594360351562, 16.00972175598145),
    @(10.6145133972168, -29.43916702270508, -53.48576354980469),
    @(-6.689743518829346, -5.009637832641602, -2.467369556427002),
    @(-15.53036594390869, -35.93496704101562, 9.646455764770508),
    @(-8.54952621459961, -4.356056213378906, -4.340849876403809),
    @(0.2473421096801757, -19.03062438964844, 14.19175815582275),
    @(-66.0440902709961, -40.48963928222656, -45.51393508911133)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
